# add value precio in pedido_producto and implement this function in frontend and backend
#
# Sheet "Inventario" (sheet1):
#   - row4 "ÚLTIMA ACTUALIZACIÓN" / "28/10/2024"  ->  "STOCK ACTUAL" / 173 (number)
#   - new row5: "FECHA DE ACTUALIZACIÓN" / "2024-10-28T19:46:07.904Z"
#
# Sheet "Productos" (sheet2):
#   - re-exported product list (one stale row "holabro2dfsad3" / "Vino" removed,
#     remaining rows reordered) - rewrite the whole table.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Inventario")
$ws2 = $wb.Worksheets.Item("Productos")

# ---- Inventario sheet ----
$ws1.Range("A4").Value = "STOCK ACTUAL"
$ws1.Range("B4").Value = 173

# New row 5 - copy row 3's banded style (white/gray alternation) down first,
# then overwrite with the real values.
$ws1.Range("A3:B3").Copy($ws1.Range("A5:B5"))
$ws1.Range("A5").Value = "FECHA DE ACTUALIZACIÓN"
$ws1.Range("B5").Value = "2024-10-28T19:46:07.904Z"

# ---- Productos sheet ----
$headers = @("NOMBRE DEL PRODUCTO", "MARCA", "DESCRIPCIÓN", "CATEGORÍA", "TIPO", "CANTIDAD")

$rows = @(
    ,@("NuevoTest25", "NO REGISTRADO", "holabrofdssdtest", "Cigarrillo", "Otro", 25)
    ,@("Cerveza Artesanal Actualizada testasdasd", "NO REGISTRADO", "holabrofdssdtest", "Néctar", "Otro", 25)
    ,@("Cerveza Artesanal Actualizada tesasdasdt", "NO REGISTRADO", "testttasdasd", "Snack", "Otro", 23)
    ,@("nuevo test 2", "NO REGISTRADO", "testttasdasd", "Agua mineral", "Sin Alcohol", 20)
    ,@("Ron 2 test", "NO REGISTRADO", "asasfafasfsafsfasasffasasfda", "Ron", "Alcohólico", 5)
    ,@("Vodka", "NO REGISTRADO", "Nuevo Vodka Vegano", "Otro", "Alcohólico", 22)
    ,@("Cerveza Artesanal Actualizada test", "NO REGISTRADO", "holabrofdssdtest", "Cerveza", "Sin Alcohol", 50)
    ,@("Ron", "testqwrwqd", "Ron Vodka cualificado por admins", "Ron", "Alcohólico", 3)
)

for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws2.Cells.Item(1, $c + 1).Value = $headers[$c]
}

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# The previous export had 10 data rows (row 10), the new export only has 9 -
# delete the now-unused trailing row so the sheet's used range shrinks back
# to match.
$ws2.Rows.Item(10).Delete()
